$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.7
$ws.Range("I2").Value = 2.63
$ws.Range("J2").Value = 3.25
$ws.Range("L2").Value = 3.25
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 7.5
$ws.Range("Z2").Value = 12
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 23
$ws.Range("AD2").Value = 34
$ws.Range("AG2").Value = 15
$ws.Range("AI2").Value = 800
$ws.Range("AJ2").Value = 8
$ws.Range("AK2").Value = 12
$ws.Range("AL2").Value = 11
$ws.Range("AM2").Value = 26
$ws.Range("AN2").Value = 23

# Row 3 updates
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6
$ws.Range("AI3").Value = 1000
$ws.Range("AP3").Value = 1.98
$ws.Range("AQ3").Value = 1.88

# Row 4 updates
$ws.Range("G4").Value = 2.63
$ws.Range("H4").Value = 3
$ws.Range("J4").Value = 3.25
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 11
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.75
$ws.Range("Q4").Value = 1.95
$ws.Range("R4").Value = 1.9
$ws.Range("S4").Value = 3.25
$ws.Range("T4").Value = 1.33
$ws.Range("W4").Value = 1.73
$ws.Range("X4").Value = 2
$ws.Range("AE4").Value = 9.5
$ws.Range("AK4").Value = 13
$ws.Range("AM4").Value = 26
